# Adds daily COVID-19 boletim rows for 2021-05-22 through 2021-05-30
# (date serials 44338-44346) to the bottom of "Planilha1", matching the
# commit "add data until May 30th".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Each inner array is one new row: Date, B..J daily/cumulative counters,
# in the same column order as the existing data (A:J).
$newRows = @(
    @(44338, 10437, 180, 4476, 15093, 3815, 558, 25, 533, 103),
    @(44339, 10473, 122, 4502, 15097, 3857, 542, 24, 518, 103),
    @(44340, 10522, 133, 4541, 15196, 3928, 509, 24, 485, 104),
    @(44341, 10565, 184, 4581, 15330, 3989, 488, 24, 464, 104),
    @(44342, 10648, 220, 4655, 15523, 4036, 515, 24, 491, 104),
    @(44343, 10773, 188, 4751, 15712, 4059, 588, 27, 561, 104),
    @(44344, 10826, 221, 4805, 15852, 4107, 592, 25, 567, 106),
    @(44345, 10934, 174, 4904, 16012, 4182, 614, 26, 588, 108),
    @(44346, 10934, 172, 4906, 16012, 4213, 584, 27, 557, 109)
)

$startRow = 325
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value2 = $row[$c]
    }
}

$lastRow = $startRow + $newRows.Length - 1

# Scroll the frozen view down to the newly added rows and select the
# next empty row below the new data, mirroring what the author saw
# after pasting/typing the new rows in Excel.
$ws.Select()
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("A2").Select()
$win.FreezePanes = $true
$win.ScrollRow = $startRow

$selRow = $lastRow + 1
$ws.Range("A" + $selRow).Select()
